# "updated template of students"
# - Rename header "Классы" -> "Класс"
# - Birth-date column (H) switches from Excel date-serial numbers to plain text dates
# - Selection moves to H5
# - A few row heights shrink slightly (content reflow after the format change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header rename: "Классы" -> "Класс" (column E, header row 2)
$ws.Range("E2").Value = "Класс"

# 2) Birth dates (column H) become literal text instead of Excel dates.
#    Force Text number format first so the literal strings are not
#    re-interpreted as date serials.
$ws.Range("H3:H6").NumberFormat = "@"
$ws.Range("H3").Value = "11.02.1996"
$ws.Range("H4").Value = "05.10.1995"
$ws.Range("H5").Value = "08.11.1995"
$ws.Range("H6").Value = "25.09.1996"

# 3) Row heights tighten up now that the date column no longer wraps.
$ws.Rows.Item(3).RowHeight = 14.15
$ws.Rows.Item(4).RowHeight = 14.15
$ws.Rows.Item(5).RowHeight = 14.15
$ws.Rows.Item(6).RowHeight = 23.85

# 4) Column widths adjust slightly.
$ws.Columns.Item(2).ColumnWidth = 12.17
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(7).ColumnWidth = 14.67
$ws.Columns.Item(8).ColumnWidth = 21.33
$ws.Columns.Item(9).ColumnWidth = 34.33
$ws.Columns.Item(10).ColumnWidth = 15.5
$ws.Columns.Item(11).ColumnWidth = 20.17
$ws.Columns.Item(13).ColumnWidth = 15.5
$ws.Columns.Item(14).ColumnWidth = 14
$ws.Columns.Item(15).ColumnWidth = 16.5
$ws.Columns.Item(17).ColumnWidth = 16.5

# 5) Selection ends on H5.
$ws.Range("H5").Select() | Out-Null
